$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 23-26 (the emoji/symbol rows), which shifts rows 27-33
# up to become rows 23-29, matching the target diff.
$ws.Range("A23:B26").EntireRow.Delete()
